# zve test correction #7
# Applies the corrections described in the commit:
#  - Fix off-by-one in column X (MAX(W-AB,0) must use the SAME row's AB, not the next row's)
#  - Fix the "1908" typo -> "1308" in U14/U16 (matches the pattern already used in U18)
#  - Fix AA16 threshold 1000 -> 920 (matches the pattern already used in AA18/AA8:AA13)
#  - Fix AA18 to use MAX(...) instead of MIN(...) for the employment-income add-back, and
#    parenthesize the (12*C18) term
#  - Add a new helper cell AA28 that keeps the old MIN(...) computation that used to live
#    inside AA18, so it is not lost
#  - Fill in the missing P8 cell (it was skipped, unlike every other row)
#  - Update the active selection/view to reflect where the author was last looking

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- missing data point -------------------------------------------------
$ws.Range("P8").Value = 0

# --- column X: MAX(W{r}-AB{r},0) instead of MAX(W{r}-AB{r+1},0) --------
$ws.Range("X3").Formula = "=MAX(W3-AB3,0)"
$ws.Range("X4:X25").Formula = "=MAX(W4-AB4,0)"

# --- U14 / U16: 1908 -> 1308 --------------------------------------------
$ws.Range("U14").Formula = "=MAX(AA14-Z14-36-1308,0)"
$ws.Range("U16").Formula = "=MAX(AA16-Z16-36-1308,0)"

# --- AA16: 1000 -> 920 ----------------------------------------------------
$ws.Range("AA16").Formula = "=(12*D16+(12*C16 - 920)*(C16>450))"

# --- AA18: MIN -> MAX, add parens around (12*C18) ------------------------
$ws.Range("AA18").Formula = "=(12*D18+((12*C18)-920)*(C18>450)+MAX((12*E18)-51-1370,0))"

# --- new row 28: keep the old MIN(...) expression that used to sit in AA18
$ws.Range("AA28").Formula = "=MIN((12*E18)-51-1370,0)"

# --- update view/selection to match the author's last position ----------
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("V18").Select()
